$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 144, shifting existing rows 144..276 down to 145..277.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly sample.
$ws.Cells.Item(144, 1).Value = 8
$ws.Cells.Item(144, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44589
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 100114013
$ws.Cells.Item(144, 7).Value = "Zanahoria"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 800
$ws.Cells.Item(144, 11).Value = 5500
$ws.Cells.Item(144, 12).Value = 6000
$ws.Cells.Item(144, 13).Value = 5750
$ws.Cells.Item(144, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(144, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(144, 16).Value = 288
$ws.Cells.Item(144, 17).Value = 20
$ws.Cells.Item(144, 18).Value = "Hortaliza"
